$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.254.12'
$ws.Range('E2').Value = '  -4.84%  '
$ws.Range('D3').Value = '3.243.22'
$ws.Range('E3').Value = '  -7.96%  '
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '595.75'
$ws.Range('E5').Value = '  -1.67%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '153.02'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '3.230.83'
$ws.Range('E8').Value = '  -8.21%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.545'
$ws.Range('E9').Value = '  -11.16%  '
$ws.Range('E10').Value = '  -11.56%  '
$ws.Range('E11').Value = '  -10.18%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.493'
$ws.Range('E12').Value = '  -15.96%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '39.10'
$ws.Range('E13').Value = '  -15.82%  '
$ws.Range('E14').Value = '  -12.41%  '
$ws.Range('D15').Value = '3.768.42'
$ws.Range('E15').Value = '  -8.10%  '
$ws.Range('D16').Value = '67.250.56'
$ws.Range('E16').Value = '  -4.93%  '
$ws.Range('D17').Value = '3.241.04'
$ws.Range('E17').Value = '  -7.93%  '
$ws.Range('E18').Value = '  -4.56%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '533.70'
$ws.Range('E19').Value = '  -12.72%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.09'
$ws.Range('E20').Value = '  -14.82%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.91'
$ws.Range('E21').Value = '  -15.21%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.759'
$ws.Range('E22').Value = '  -13.85%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.93'
$ws.Range('E23').Value = '  -11.95%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.78'
$ws.Range('E24').Value = '  -11.95%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '85.66'
$ws.Range('E25').Value = '  -13.54%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.998'
$ws.Range('E26').Value = '  -0.22%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.20'
$ws.Range('E27').Value = '  -14.08%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.20'
$ws.Range('E28').Value = '  -14.67%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.02'
$ws.Range('E29').Value = '  -11.46%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '29.15'
$ws.Range('E30').Value = '  -14.41%  '
$ws.Range('E31').Value = '  -10.16%  '
$ws.Range('E32').Value = '  -9.57%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '532.06'
$ws.Range('E33').Value = '  -17.06%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.76'
$ws.Range('E34').Value = '  -15.80%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.44'
$ws.Range('E35').Value = '  -20.34%  '
$ws.Range('E36').Value = '  +0.30%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '53.50'
$ws.Range('E37').Value = '  -5.90%  '
$ws.Range('E38').Value = '  -10.37%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0860'
$ws.Range('E39').Value = '  -13.73%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '9.30'
$ws.Range('E40').Value = '  -13.72%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.123'
$ws.Range('E41').Value = '  -13.44%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.76'
$ws.Range('E42').Value = '  -22.98%  '
$ws.Range('D43').Value = '2.939.39'
$ws.Range('E43').Value = '  -12.90%  '
$ws.Range('E44').Value = '  -15.03%  '
$ws.Range('D45').Value = '0.0₃0589'
$ws.Range('E45').Value = '  -20.89%  '
$ws.Range('E46').Value = '  -15.65%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '26.50'
$ws.Range('E47').Value = '  -17.71%  '
$ws.Range('E48').Value = '  -16.22%  '
$ws.Range('E49').Value = '  -0.08%  '
$ws.Range('E50').Value = '  -11.83%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '118.70'
$ws.Range('E51').Value = '  -11.16%  '
